# [LC-850] Release of LetsCo OS v1.3.0
#
# The KPI "name" codes in column B were renamed from a 1-digit suffix
# scheme to a zero-padded 2-digit suffix scheme:
#   GP1 -> GP01, GP2 -> GP02, GP3 -> GP03
#   BP1 -> BP01, BP2 -> BP02, BP3 -> BP03
#
# The sheet's saved view/selection is also reset back to cell A1 (it was
# previously left scrolled/selected at G1 / I4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$replacements = @{
    "GP1" = "GP01"
    "GP2" = "GP02"
    "GP3" = "GP03"
    "BP1" = "BP01"
    "BP2" = "BP02"
    "BP3" = "BP03"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row  # xlUp
if ($lastRow -lt 1) { $lastRow = 1 }
if ($lastRow -gt 1000) { $lastRow = 1000 }  # sane upper bound safety net

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    $text = $cell.Text
    if ($replacements.ContainsKey($text)) {
        $cell.Value = $replacements[$text]
    }
}

# Reset the saved view/selection back to the top-left corner (A1) instead
# of the previously scrolled/selected G1 (topLeftCell) / I4 (selection) state.
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
